$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows of accelerometer samples above the existing data (old row 2),
# shifting the existing data down by 8 rows.
$ws.Range("A2:C9").Insert(-4121)
$ws.Range("A2:C9").ClearFormats()

# Fill in the values for the newly inserted rows.
$ws.Range("A2").Value = -0.5779368877410893
$ws.Range("B2").Value = 1.070879459381104
$ws.Range("C2").Value = 0.1698004633188247
$ws.Range("A3").Value = -0.6250030517578123
$ws.Range("B3").Value = 1.073733139038086
$ws.Range("C3").Value = 0.1257202506065367
$ws.Range("A4").Value = -0.5019012451171875
$ws.Range("B4").Value = 1.114973473548889
$ws.Range("C4").Value = 0.08085805475711817
$ws.Range("A5").Value = -0.5343909263610841
$ws.Range("B5").Value = 1.139204859733582
$ws.Range("C5").Value = 0.1443376690149308
$ws.Range("A6").Value = -0.5579452037811278
$ws.Range("B6").Value = 1.112600553035736
$ws.Range("C6").Value = 0.2124309107661247
$ws.Range("A7").Value = -0.4796955108642578
$ws.Range("B7").Value = 1.016827774047851
$ws.Range("C7").Value = 0.1028751075267787
$ws.Range("A8").Value = -0.5379581451416018
$ws.Range("B8").Value = 0.9855325698852542
$ws.Range("C8").Value = -0.2731702357530603
$ws.Range("A9").Value = -0.6476110458374021
$ws.Range("B9").Value = 1.080279231071473
$ws.Range("C9").Value = -0.8854551434516924

# Append 2 more new rows of samples at the end of the data.
$ws.Range("A30").Value = 0.00381779670715305
$ws.Range("B30").Value = 1.210070580244063
$ws.Range("C30").Value = -0.2140652965754264
$ws.Range("A31").Value = -0.03790302276611289
$ws.Range("B31").Value = 1.083934617042542
$ws.Range("C31").Value = -0.04758519232273038
